{"js": "// Appends the Q4 \"hierarchical clustering\" R-markdown block right after the\n// existing final paragraph of the document (\"table(clusters_h3)\"), matching\n// the diff that inserts 30 new paragraphs before the closing sectPr.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Anchor on the current last paragraph in the document\n// (the existing \"table(clusters_h3)\" line).\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\nconst newParagraphText = [\n  \"\",\n  \"Q4\",\n  \"# FIRST WAY I DID IT BEFORE THE ALTERNATIVE WAY\",\n  \"#Q4 Task: Conduct a hierarchical clustering analysis. Be sure to specify the linkage method used. Within your analysis, make sure you do both of the following:\",\n  \"#  1. Determine the optimal number of clusters using a clear, data-driven strategy.\",\n  \"#  2. Describe the composition of each cluster in terms of the original input features\",\n  \"\",\n  \"data(USArrests)\",\n  \"\",\n  \"# NTS: First I need to create Dissimilarity matrix\",\n  \"diss_matrix <- dist(USArrests, method = \\\"euclidean\\\")\",\n  \"\",\n  \"#NTS: Then hierarchical clustering using Complete Linkage\",\n  \"clusters_h <- hclust(diss_matrix, method = \\\"complete\\\" )\",\n  \"\",\n  \"# Plot the obtained dendrogram\",\n  \"plot(clusters_h, cex = 0.4, hang = -1)\",\n  \"\",\n  \"#NTS: Now I need to choose where to cut across my dendrogram to choose my number of clusters. In order to do that I have to  creating a function to use within clusGap. I am now using an average linkage. After running the function, I will plot the gapstat to visualize it\",\n  \"\",\n  \"hclusCut <- function(x, k) list(cluster = cutree(hclust(dist(x, method = \\\"euclidian\\\"), method = \\\"complete\\\"), k = k))\",\n  \"\",\n  \"gap_stat <- clusGap(USArrests, FUN = hclusCut, K.max = 10, B = 50)\",\n  \"fviz_gap_stat(gap_stat)\",\n  \"\",\n  \"#NTS: Now, use the number of clusters from gap statistic to obtain cluster assignment for each observation\",\n  \"clusters_h3 = cutree(clusters_h, k = 3)\",\n  \"table(clusters_h3)\",\n  \"\",\n  \"#Now I want to take those values and put it back onto my original dataset\",\n];\n\nfor (const text of newParagraphText) {\n  anchor = anchor.insertParagraph(text, Word.InsertLocation.after);\n}\n\nawait context.sync();\n", "ps1": "# Appends the Q4 \"hierarchical clustering\" R-markdown block right after the\n# existing final paragraph of the document (\"table(clusters_h3)\"), matching\n# the diff that inserts 30 new paragraphs before the closing sectPr.\n\n$d = $word.ActiveDocument\n\n$newParagraphText = @(\n    '',\n    'Q4',\n    '# FIRST WAY I DID IT BEFORE THE ALTERNATIVE WAY',\n    '#Q4 Task: Conduct a hierarchical clustering analysis. Be sure to specify the linkage method used. Within your analysis, make sure you do both of the following:',\n    '#  1. Determine the optimal number of clusters using a clear, data-driven strategy.',\n    '#  2. Describe the composition of each cluster in terms of the original input features',\n    '',\n    'data(USArrests)',\n    '',\n    '# NTS: First I need to create Dissimilarity matrix',\n    'diss_matrix <- dist(USArrests, method = \"euclidean\")',\n    '',\n    '#NTS: Then hierarchical clustering using Complete Linkage',\n    'clusters_h <- hclust(diss_matrix, method = \"complete\" )',\n    '',\n    '# Plot the obtained dendrogram',\n    'plot(clusters_h, cex = 0.4, hang = -1)',\n    '',\n    '#NTS: Now I need to choose where to cut across my dendrogram to choose my number of clusters. In order to do that I have to  creating a function to use within clusGap. I am now using an average linkage. After running the function, I will plot the gapstat to visualize it',\n    '',\n    'hclusCut <- function(x, k) list(cluster = cutree(hclust(dist(x, method = \"euclidian\"), method = \"complete\"), k = k))',\n    '',\n    'gap_stat <- clusGap(USArrests, FUN = hclusCut, K.max = 10, B = 50)',\n    'fviz_gap_stat(gap_stat)',\n    '',\n    '#NTS: Now, use the number of clusters from gap statistic to obtain cluster assignment for each observation',\n    'clusters_h3 = cutree(clusters_h, k = 3)',\n    'table(clusters_h3)',\n    '',\n    '#Now I want to take those values and put it back onto my original dataset'\n)\n\n# Anchor on the current last paragraph in the document\n# (the existing \"table(clusters_h3)\" line).\n$anchor = $d.Paragraphs.Last.Range\n\nforeach ($t in $newParagraphText) {\n    $anchor.InsertParagraphAfter()\n    $anchor = $d.Paragraphs.Last.Range\n    if ($t -ne '') {\n        $anchor.Text = $t\n    }\n}\n"}
